$d = $word.ActiveDocument

# 1) Remove the empty paragraph that immediately follows the
#    "Voi tat ca vai tro xem thong tin..." paragraph (it sits between that
#    paragraph and the "Buoc 1: O thanh Menu..." paragraph).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text.TrimEnd([char]13, [char]7) -match "liên lạc") {
    $targetIndex = $i + 1
    break
  }
}
if ($targetIndex -gt 0) {
  $empty = $d.Paragraphs.Item($targetIndex)
  if ($empty.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
    $empty.Range.Delete()
  }
}

# 2) Mark the runs holding the two trailing screenshots as NoProof
#    (adds <w:rPr><w:noProof/></w:rPr> to those runs), matching them by
#    their picture size so paragraph-index shifts from step 1 don't matter.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
  $shape = $d.InlineShapes.Item($i)
  $height = [math]::Round($shape.Height, 1)
  if ($height -eq 110.7 -or $height -eq 88.6) {
    $shape.Range.NoProofing = $true
  }
}
